# Trim trailing whitespace from company names in column C across all sheets.
# The source data had values like "CORPORACIÓN ENRIQUE & JANIS " and
# "MULTITRANS RR " with a trailing space (forcing xml:space="preserve").
# This normalizes them by removing the trailing whitespace, as part of
# unifying the unit JSON with the rest of the program.

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $usedRange = $ws.UsedRange
    $lastRow = $usedRange.Rows.Count

    for ($r = 1; $r -le $lastRow; $r++) {
        $cell = $ws.Cells.Item($r, 3)  # Column C
        $val = $cell.Value2
        if ($null -ne $val -and $val -is [string]) {
            $trimmed = $val.TrimEnd()
            if ($trimmed -ne $val) {
                $cell.Value2 = $trimmed
            }
        }
    }
}
